$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 383.63635
$ws.Range("I2").Value = 113.333336
$ws.Range("J2").Value = 485
$ws.Range("K2").Value = 113.333336
$ws.Range("L2").Value = 485
$ws.Range("M2").Value = -0.3333360000000027
$ws.Range("N2").Value = -711

$ws.Range("H64").Value = 4127.7144
$ws.Range("J64").Value = 3974
$ws.Range("L64").Value = 3974
$ws.Range("N64").Value = -4470

$ws.Range("H67").Value = 4127.7144
$ws.Range("J67").Value = 3974
$ws.Range("L67").Value = 3974
$ws.Range("N67").Value = -5690

$ws.Range("H76").Value = 9334.111000000001
$ws.Range("I76").Value = 8000.6
$ws.Range("K76").Value = 8000.6
$ws.Range("M76").Value = -7685.6

$ws.Range("H79").Value = 9334.111000000001
$ws.Range("I79").Value = 8000.6
$ws.Range("K79").Value = 8000.6
$ws.Range("M79").Value = -6908.6

$ws.Range("H92").Value = 696.8261
$ws.Range("I92").Value = 576.35
$ws.Range("K92").Value = 576.35
$ws.Range("M92").Value = 671.65

$ws.Range("H111").Value = 1428.875
$ws.Range("I111").Value = 1316.5
$ws.Range("K111").Value = 3949.5
$ws.Range("M111").Value = -882.5

$ws.Range("H135").Value = 34483456
$ws.Range("I135").Value = 341
$ws.Range("J135").Value = 166668740
$ws.Range("K135").Value = 3069
$ws.Range("L135").Value = 1500018660
$ws.Range("M135").Value = -534
$ws.Range("N135").Value = -1500023730

$ws.Range("H138").Value = 449851.5
$ws.Range("I138").Value = 1181.4445
$ws.Range("J138").Value = 714640.4
$ws.Range("K138").Value = 3544.3335
$ws.Range("L138").Value = 2143921.2
$ws.Range("M138").Value = 1595.6665
$ws.Range("N138").Value = -2154201.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3293.7778
$ws.Range("I2").Value = 501
$ws.Range("J2").Value = 9641
$ws.Range("K2").Value = 501
$ws.Range("L2").Value = 9641
$ws.Range("M2").Value = -388
$ws.Range("N2").Value = -9867

$ws.Range("H32").Value = 3968.4375
$ws.Range("I32").Value = 3606.3035
$ws.Range("J32").Value = 6503.375
$ws.Range("K32").Value = 3606.3035
$ws.Range("L32").Value = 6503.375
$ws.Range("M32").Value = -3319.3035
$ws.Range("N32").Value = -7077.375

$ws.Range("H45").Value = 1102.4584
$ws.Range("I45").Value = 1057.1875
$ws.Range("K45").Value = 1057.1875
$ws.Range("M45").Value = -680.1875

$ws.Range("H61").Value = 1383.2222
$ws.Range("I61").Value = 1273.88
$ws.Range("J61").Value = 2750
$ws.Range("K61").Value = 1273.88
$ws.Range("L61").Value = 2750
$ws.Range("M61").Value = -1061.88
$ws.Range("N61").Value = -3174

$ws.Range("H116").Value = 3293.7778
$ws.Range("I116").Value = 501
$ws.Range("J116").Value = 9641
$ws.Range("K116").Value = 501
$ws.Range("L116").Value = 9641
$ws.Range("M116").Value = 1793
$ws.Range("N116").Value = -14229

$ws.Range("H123").Value = 48601
$ws.Range("J123").Value = 48601
$ws.Range("L123").Value = 48601
$ws.Range("N123").Value = -58401

$ws.Range("H132").Value = 1088.0892
$ws.Range("I132").Value = 834
$ws.Range("J132").Value = 1928.5385
$ws.Range("K132").Value = 2502
$ws.Range("L132").Value = 5785.6155
$ws.Range("M132").Value = 28
$ws.Range("N132").Value = -10845.6155

$ws.Range("H136").Value = 1383.2222
$ws.Range("I136").Value = 1273.88
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 3821.64
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -1271.64
$ws.Range("N136").Value = -13350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3293.7778
$ws.Range("I3").Value = 501
$ws.Range("J3").Value = 9641
$ws.Range("K3").Value = 501
$ws.Range("L3").Value = 9641
$ws.Range("M3").Value = -387
$ws.Range("N3").Value = -9869

$ws.Range("H20").Value = 1452.9565
$ws.Range("I20").Value = 1018.2778
$ws.Range("K20").Value = 1018.2778
$ws.Range("M20").Value = -771.2778

$ws.Range("H86").Value = 3459.6
$ws.Range("I86").Value = 3731.0527
$ws.Range("K86").Value = 3731.0527
$ws.Range("M86").Value = -2608.0527

$ws.Range("H89").Value = 3459.6
$ws.Range("I89").Value = 3731.0527
$ws.Range("K89").Value = 18655.2635
$ws.Range("M89").Value = -13039.2635

$ws.Range("H94").Value = 25000854
$ws.Range("I94").Value = 35714816
$ws.Range("J94").Value = 1616.6666
$ws.Range("K94").Value = 35714816
$ws.Range("L94").Value = 1616.6666
$ws.Range("M94").Value = -35714365
$ws.Range("N94").Value = -2518.6666

$ws.Range("H134").Value = 3879.1191
$ws.Range("I134").Value = 1029.4062
$ws.Range("J134").Value = 12998.2
$ws.Range("K134").Value = 3088.2186
$ws.Range("L134").Value = 38994.60000000001
$ws.Range("M134").Value = -553.2185999999997
$ws.Range("N134").Value = -44064.60000000001

$ws.Range("H135").Value = 38411.555
$ws.Range("J135").Value = 38411.555
$ws.Range("L135").Value = 38411.555
$ws.Range("N135").Value = -48551.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1266.8334
$ws.Range("I31").Value = 1190.6792
$ws.Range("J31").Value = 1577.3077
$ws.Range("K31").Value = 1190.6792
$ws.Range("L31").Value = 1577.3077
$ws.Range("M31").Value = -895.6792
$ws.Range("N31").Value = -2167.3077

$ws.Range("H34").Value = 1266.8334
$ws.Range("I34").Value = 1190.6792
$ws.Range("J34").Value = 1577.3077
$ws.Range("K34").Value = 1190.6792
$ws.Range("L34").Value = 1577.3077
$ws.Range("M34").Value = -988.6792
$ws.Range("N34").Value = -1981.3077

$ws.Range("H58").Value = 703
$ws.Range("I58").Value = 614.9429
$ws.Range("J58").Value = 1216.6666
$ws.Range("K58").Value = 614.9429
$ws.Range("L58").Value = 1216.6666
$ws.Range("M58").Value = -411.9429
$ws.Range("N58").Value = -1622.6666

$ws.Range("H132").Value = 4868.515
$ws.Range("I132").Value = 5067.148
$ws.Range("K132").Value = 15201.444
$ws.Range("M132").Value = -12671.444

$ws.Range("H134").Value = 1481.3529
$ws.Range("I134").Value = 1528.48
$ws.Range("J134").Value = 1350.4445
$ws.Range("K134").Value = 4585.440000000001
$ws.Range("L134").Value = 4051.3335
$ws.Range("M134").Value = -2050.440000000001
$ws.Range("N134").Value = -9121.333500000001

$ws.Range("H136").Value = 703
$ws.Range("I136").Value = 614.9429
$ws.Range("J136").Value = 1216.6666
$ws.Range("K136").Value = 1844.8287
$ws.Range("L136").Value = 3649.9998
$ws.Range("M136").Value = 705.1713
$ws.Range("N136").Value = -8749.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1229.5667
$ws.Range("I5").Value = 1262.2963
$ws.Range("K5").Value = 3786.8889
$ws.Range("M5").Value = -3674.8889

$ws.Range("H32").Value = 3200
$ws.Range("I32").Value = 3200
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 9600
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9317
$ws.Range("N32").ClearContents()

$ws.Range("H38").Value = 55.76923
$ws.Range("J38").Value = 79.333336
$ws.Range("L38").Value = 238.000008
$ws.Range("N38").Value = -932.000008

$ws.Range("H56").Value = 6394.6523
$ws.Range("I56").Value = 6394.6523
$ws.Range("K56").Value = 6394.6523
$ws.Range("M56").Value = -5864.6523

$ws.Range("H131").Value = 12822835
$ws.Range("J131").Value = 2352.4934
$ws.Range("L131").Value = 7057.4802
$ws.Range("N131").Value = -17137.4802

$ws.Range("H135").Value = 1229.5667
$ws.Range("I135").Value = 1262.2963
$ws.Range("K135").Value = 11360.6667
$ws.Range("M135").Value = -8825.6667

$ws.Range("H136").Value = 1617.2778
$ws.Range("J136").Value = 3326.2
$ws.Range("L136").Value = 9978.599999999999
$ws.Range("N136").Value = -20178.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H126").Value = 2204.75
$ws.Range("I126").Value = 1763.1818
$ws.Range("J126").Value = 2744.4443
$ws.Range("K126").Value = 5289.5454
$ws.Range("L126").Value = 8233.332900000001
$ws.Range("M126").Value = -2819.5454
$ws.Range("N126").Value = -13173.3329

$ws.Range("H128").Value = 38000
$ws.Range("J128").Value = 38000
$ws.Range("L128").Value = 38000
$ws.Range("N128").Value = -47960

$ws.Range("H132").Value = 1797.0555
$ws.Range("I132").Value = 1357.4103
$ws.Range("J132").Value = 2940.1333
$ws.Range("K132").Value = 4072.2309
$ws.Range("L132").Value = 8820.3999
$ws.Range("M132").Value = -1542.2309
$ws.Range("N132").Value = -13880.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1422.5
$ws.Range("J22").Value = 1422.5
$ws.Range("L22").Value = 1422.5
$ws.Range("N22").Value = -2012.5

$ws.Range("H27").Value = 1422.5
$ws.Range("J27").Value = 1422.5
$ws.Range("L27").Value = 1422.5
$ws.Range("N27").Value = -1636.5

$ws.Range("H46").Value = 6928.5713
$ws.Range("I46").Value = 2600
$ws.Range("K46").Value = 2600
$ws.Range("M46").Value = -2412

$ws.Range("H132").Value = 20952.059
$ws.Range("I132").Value = 1262.625
$ws.Range("J132").Value = 52455.15
$ws.Range("K132").Value = 3787.875
$ws.Range("L132").Value = 157365.45
$ws.Range("M132").Value = -1257.875
$ws.Range("N132").Value = -162425.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 418.73685
$ws.Range("J107").Value = 452.57144
$ws.Range("L107").Value = 1357.71432
$ws.Range("N107").Value = -5197.71432

$ws.Range("H132").Value = 2548.8445
$ws.Range("I132").Value = 2434.6
$ws.Range("J132").Value = 2948.7
$ws.Range("K132").Value = 7303.799999999999
$ws.Range("L132").Value = 8846.099999999999
$ws.Range("M132").Value = -4773.799999999999
$ws.Range("N132").Value = -13906.1

$ws.Range("H136").Value = 833
$ws.Range("I136").Value = 558.5
$ws.Range("J136").Value = 1748
$ws.Range("K136").Value = 1875.5
$ws.Range("L136").Value = 5244
$ws.Range("M136").Value = 874.5
$ws.Range("N136").Value = -10344
